$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new marker columns in J1 and K1
$ws.Range("J1").Value = "marker_1"
$ws.Range("K1").Value = "marker_2"

# Update selection to reflect the newly added columns
$ws.Range("J1:K1").Select()
